$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 243.86667
$ws.Range("I33").Value = 225.57143
$ws.Range("K33").Value = 225.57143
$ws.Range("M33").Value = 3.428570000000008

$ws.Range("H69").Value = 10199.8
$ws.Range("I69").Value = 9000
$ws.Range("K69").Value = 27000
$ws.Range("M69").Value = -26126

$ws.Range("H72").Value = 10199.8
$ws.Range("I72").Value = 9000
$ws.Range("K72").Value = 81000
$ws.Range("M72").Value = -76632

$ws.Range("H100").Value = 3170.8
$ws.Range("I100").Value = 923.8333
$ws.Range("K100").Value = 923.8333
$ws.Range("M100").Value = -382.8333

$ws.Range("H137").Value = 4155.4146
$ws.Range("I137").Value = 1739.5278
$ws.Range("K137").Value = 5218.5834
$ws.Range("M137").Value = -2668.5834

$ws.Range("H138").Value = 3372.889
$ws.Range("I138").Value = 2613.4
$ws.Range("K138").Value = 7840.200000000001
$ws.Range("M138").Value = -2700.200000000001

$ws.Range("H141").Value = 5821.2383
$ws.Range("I141").Value = 2231.9285
$ws.Range("K141").Value = 6695.7855
$ws.Range("M141").Value = -1515.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 736.5599999999999
$ws.Range("J2").Value = 737
$ws.Range("L2").Value = 737
$ws.Range("N2").Value = -963

$ws.Range("H32").Value = 218667.72
$ws.Range("I32").Value = 226075.17
$ws.Range("K32").Value = 226075.17
$ws.Range("M32").Value = -225788.17

$ws.Range("H61").Value = 3114.2354
$ws.Range("I61").Value = 3121.375
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 3121.375
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -2909.375
$ws.Range("N61").Value = -3424

$ws.Range("H74").Value = 6615.4165
$ws.Range("I74").Value = 3999.9023
$ws.Range("K74").Value = 3999.9023
$ws.Range("M74").Value = -3125.9023

$ws.Range("H77").Value = 6615.4165
$ws.Range("I77").Value = 3999.9023
$ws.Range("K77").Value = 19999.5115
$ws.Range("M77").Value = -15631.5115

$ws.Range("H102").Value = 1681.1428
$ws.Range("I102").Value = 1295.5834
$ws.Range("K102").Value = 1295.5834
$ws.Range("M102").Value = 326.4166

$ws.Range("H116").Value = 736.5599999999999
$ws.Range("J116").Value = 737
$ws.Range("L116").Value = 737
$ws.Range("N116").Value = -5325

$ws.Range("H132").Value = 4411.1064
$ws.Range("I132").Value = 2627.4119
$ws.Range("K132").Value = 7882.2357
$ws.Range("M132").Value = -5352.2357

$ws.Range("H136").Value = 3114.2354
$ws.Range("I136").Value = 3121.375
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 9364.125
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -6814.125
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 736.5599999999999
$ws.Range("J3").Value = 737
$ws.Range("L3").Value = 737
$ws.Range("N3").Value = -965

$ws.Range("H5").Value = 1468
$ws.Range("I5").Value = 761.6
$ws.Range("J5").Value = 5000
$ws.Range("K5").Value = 761.6
$ws.Range("L5").Value = 5000
$ws.Range("M5").Value = -648.6
$ws.Range("N5").Value = -5226

$ws.Range("H35").Value = 35123.5
$ws.Range("J35").Value = 34999.668
$ws.Range("L35").Value = 34999.668
$ws.Range("N35").Value = -35619.668

$ws.Range("H94").Value = 4983.9
$ws.Range("I94").Value = 4542.375
$ws.Range("J94").Value = 6750
$ws.Range("K94").Value = 4542.375
$ws.Range("L94").Value = 6750
$ws.Range("M94").Value = -4091.375
$ws.Range("N94").Value = -7652

$ws.Range("H99").Value = 15785.714
$ws.Range("I99").Value = 15785.714
$ws.Range("K99").Value = 15785.714
$ws.Range("M99").Value = -14287.714

$ws.Range("H105").Value = 7955.8335
$ws.Range("I105").Value = 1949.5
$ws.Range("J105").Value = 19968.5
$ws.Range("K105").Value = 1949.5
$ws.Range("L105").Value = 19968.5
$ws.Range("M105").Value = -202.5
$ws.Range("N105").Value = -23462.5

$ws.Range("H134").Value = 6738.75
$ws.Range("I134").Value = 7046
$ws.Range("K134").Value = 21138
$ws.Range("M134").Value = -18603

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 9255.5
$ws.Range("I25").Value = 9255.5
$ws.Range("K25").Value = 9255.5
$ws.Range("M25").Value = -9081.5

$ws.Range("H58").Value = 3375.425
$ws.Range("I58").Value = 2795.08
$ws.Range("K58").Value = 2795.08
$ws.Range("M58").Value = -2592.08

$ws.Range("H105").Value = 19857.143
$ws.Range("I105").Value = 19857.143
$ws.Range("K105").Value = 19857.143
$ws.Range("M105").Value = -18110.143

$ws.Range("H136").Value = 3375.425
$ws.Range("I136").Value = 2795.08
$ws.Range("K136").Value = 8385.24
$ws.Range("M136").Value = -5835.24

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 44642.145
$ws.Range("J37").Value = 44642.145
$ws.Range("L37").Value = 133926.435
$ws.Range("N37").Value = -134150.435

$ws.Range("H107").Value = 4696.75
$ws.Range("J107").Value = 5059.3
$ws.Range("L107").Value = 15177.9
$ws.Range("N107").Value = -19017.9

$ws.Range("H134").Value = 8117.75
$ws.Range("I134").Value = 4411
$ws.Range("J134").Value = 12883.571
$ws.Range("K134").Value = 13233
$ws.Range("L134").Value = 38650.713
$ws.Range("M134").Value = -8163
$ws.Range("N134").Value = -48790.713

$ws.Range("H137").Value = 7117.2104
$ws.Range("J137").Value = 8213.385
$ws.Range("L137").Value = 24640.155
$ws.Range("N137").Value = -34840.155

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 2000000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 2000000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 2000000
$ws.Range("N24").Value = -2000346
$ws.Range("M24").ClearContents()

$ws.Range("H80").Value = 3291
$ws.Range("J80").Value = 3800
$ws.Range("L80").Value = 3800
$ws.Range("N80").Value = -5796

$ws.Range("H83").Value = 3291
$ws.Range("J83").Value = 3800
$ws.Range("L83").Value = 19000
$ws.Range("N83").Value = -28984

$ws.Range("H97").Value = 990.5333000000001
$ws.Range("I97").Value = 964.7
$ws.Range("K97").Value = 964.7
$ws.Range("M97").Value = -468.7

$ws.Range("H107").Value = 1243.5
$ws.Range("I107").Value = 1243.5
$ws.Range("K107").Value = 1243.5
$ws.Range("M107").Value = 676.5

$ws.Range("H126").Value = 2950.5
$ws.Range("J126").Value = 3177.5
$ws.Range("L126").Value = 9532.5
$ws.Range("N126").Value = -14472.5

$ws.Range("H132").Value = 12628.8
$ws.Range("I132").Value = 12998.759
$ws.Range("K132").Value = 38996.277
$ws.Range("M132").Value = -36466.277

$ws.Range("H138").Value = 60000
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3015.8333
$ws.Range("J7").Value = 4331.6665
$ws.Range("L7").Value = 4331.6665
$ws.Range("N7").Value = -4555.6665

$ws.Range("H9").Value = 2247.3
$ws.Range("I9").Value = 339
$ws.Range("J9").Value = 6700
$ws.Range("K9").Value = 339
$ws.Range("L9").Value = 6700
$ws.Range("M9").Value = -115
$ws.Range("N9").Value = -7148

$ws.Range("H40").Value = 21840
$ws.Range("I40").Value = 26050
$ws.Range("K40").Value = 26050
$ws.Range("M40").Value = -25914

$ws.Range("H61").Value = 8782.291999999999
$ws.Range("I61").Value = 7370.2383
$ws.Range("J61").Value = 18666.666
$ws.Range("K61").Value = 7370.2383
$ws.Range("L61").Value = 18666.666
$ws.Range("M61").Value = -7168.2383
$ws.Range("N61").Value = -19070.666

$ws.Range("H113").Value = 8782.291999999999
$ws.Range("I113").Value = 7370.2383
$ws.Range("J113").Value = 18666.666
$ws.Range("K113").Value = 7370.2383
$ws.Range("L113").Value = 18666.666
$ws.Range("M113").Value = -5200.2383
$ws.Range("N113").Value = -23006.666

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122,N122").ClearContents()

$ws.Range("H126").Value = 3015.8333
$ws.Range("J126").Value = 4331.6665
$ws.Range("L126").Value = 12994.9995
$ws.Range("N126").Value = -17934.9995

$ws.Range("H132").Value = 2821.6296
$ws.Range("I132").Value = 1782.1875
$ws.Range("J132").Value = 4333.5454
$ws.Range("K132").Value = 5346.5625
$ws.Range("L132").Value = 13000.6362
$ws.Range("M132").Value = -2816.5625
$ws.Range("N132").Value = -18060.6362

$ws.Range("H136").Value = 12937.5
$ws.Range("I136").Value = 6539.8
$ws.Range("J136").Value = 19335.2
$ws.Range("K136").Value = 19619.4
$ws.Range("L136").Value = 58005.60000000001
$ws.Range("M136").Value = -17069.4
$ws.Range("N136").Value = -63105.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2755.1904
$ws.Range("I132").Value = 2093.1428
$ws.Range("J132").Value = 4079.2856
$ws.Range("K132").Value = 6279.428400000001
$ws.Range("L132").Value = 12237.8568
$ws.Range("M132").Value = -3749.428400000001
$ws.Range("N132").Value = -17297.8568

$ws.Range("H136").Value = 6969.9287
$ws.Range("I136").Value = 7429.154
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 22287.462
$ws.Range("L136").Value = 1000
$ws.Range("M136").Value = -19737.462
$ws.Range("N136").Value = -8100
